# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" / "Salario Basico" table in B16:J45 is
# refreshed from the updated database: the period column (E) is re-sorted
# ascending (previously it ran 2003 down to 1710, now it runs 1710 up to
# 2003) and the F/G figures that travel with each period are refreshed to
# match, while the "Salario Basico" column (G) is updated to the new flat
# value for every period row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period order for E16:E45 (was descending 2003 -> 1710).
$periods = @(
    "1710", "1711", "1712",
    "1801", "1802", "1803", "1804", "1805", "1806", "1807", "1808", "1809", "1810", "1811", "1812",
    "1901", "1902", "1903", "1904", "1905", "1906", "1907", "1908", "1909", "1910", "1911", "1912",
    "2001", "2002", "2003"
)

# F column ("Valor Mora") travels together with its period row, the same
# way it did before the re-sort (first 11 periods -> 29509, next 18 -> 31249,
# last period -> 30208).
$mora = @(
    29509, 29509, 29509, 29509, 29509, 29509, 29509, 29509, 29509, 29509, 29509,
    31249, 31249, 31249, 31249, 31249, 31249, 31249, 31249, 31249, 31249, 31249, 31249, 31249, 31249, 31249, 31249, 31249, 31249,
    30208
)

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $mora[$i]
    # "Salario Basico" (G) database value refreshed to 781242 for every row.
    $ws.Range("G$row").Value = 781242
}
